$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column D header from "fieldNotes" to "eventRemarks"
$ws.Range("D1").Value = "eventRemarks"

# Header row (A1:D1) loses its bold styling but keeps its border
$ws.Range("A1:D1").Font.Bold = $false

# Move the active selection to G9, matching the saved cursor position
$ws.Range("G9").Select()
